$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.010.94"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.398.81"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'485.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'154.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +17.70%  "
$ws.Range("D9").Value = "2.415.81"
$ws.Range("E9").Value = "  -3.85%  "
$ws.Range("D10").Value = "'6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.38%  "
$ws.Range("D11").Value = "'0.0993"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "2.818.66"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").Value = "56.981.66"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'20.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").Value = "2.415.54"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("D19").Value = "'4.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").Value = "'324.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'9.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'5.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'58.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").Value = "'0.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "2.512.81"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").Value = "'7.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").Value = "0.0₃0776"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'149.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "'18.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "'1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "'3.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'0.834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("D39").Value = "'34.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +7.99%  "
$ws.Range("D41").Value = "'3.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'0.594"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("D45").Value = "'268.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'0.0528"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "'4.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.36%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.865.26"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'17.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.13%  "
